$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.300.42"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.102.53"
$ws.Range("E3").Value = "  +4.13%  "
$ws.Range("E4").Value = "  +0.07%  "
$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.28"
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = "  +2.02%  "
$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.663"
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = "  +0.32%  "
$origStyle_D8 = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.92"
$ws.Range("D8").Style = $origStyle_D8
$ws.Range("E8").Value = "  +19.30%  "
$origStyle_D9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "61.80"
$ws.Range("D9").Style = $origStyle_D9
$ws.Range("E9").Value = "  +3.32%  "
$origStyle_D10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.376"
$ws.Range("D10").Style = $origStyle_D10
$ws.Range("E10").Value = "  +2.22%  "
$origStyle_D11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0746"
$ws.Range("D11").Style = $origStyle_D11
$ws.Range("E11").Value = "  +4.05%  "
$ws.Range("E12").Value = "  +7.39%  "
$origStyle_D13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.26"
$ws.Range("D13").Style = $origStyle_D13
$ws.Range("E13").Value = "  +4.49%  "
$ws.Range("D14").Value = "2.404.57"
$ws.Range("E14").Value = "  +3.88%  "
$origStyle_D15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.840"
$ws.Range("D15").Style = $origStyle_D15
$ws.Range("E15").Value = "  +3.71%  "
$ws.Range("D16").Value = "2.100.98"
$ws.Range("E16").Value = "  +3.94%  "
$origStyle_D17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.16"
$ws.Range("D17").Style = $origStyle_D17
$ws.Range("E17").Value = "  +4.87%  "
$ws.Range("D18").Value = "37.296.10"
$ws.Range("E18").Value = "  +2.05%  "
$origStyle_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "72.55"
$ws.Range("D19").Style = $origStyle_D19
$ws.Range("E19").Value = "  +1.63%  "
$origStyle_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.36"
$ws.Range("D20").Style = $origStyle_D20
$ws.Range("E20").Value = "  +10.58%  "
$ws.Range("D21").Value = "0.0₃0841"
$ws.Range("E21").Value = "  +2.48%  "
$origStyle_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "241.43"
$ws.Range("D22").Style = $origStyle_D22
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("E23").Value = "  +7.16%  "
$ws.Range("E24").Value = "  +0.15%  "
$origStyle_D25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.48"
$ws.Range("D25").Style = $origStyle_D25
$ws.Range("E25").Value = "  +1.04%  "
$origStyle_D26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.87"
$ws.Range("D26").Style = $origStyle_D26
$ws.Range("E26").Value = "  +4.66%  "
$origStyle_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("D27").Style = $origStyle_D27
$ws.Range("E27").Value = "  +7.49%  "
$origStyle_D28 = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.75"
$ws.Range("D28").Style = $origStyle_D28
$ws.Range("E28").Value = "  +5.30%  "
$origStyle_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.01"
$ws.Range("D29").Style = $origStyle_D29
$ws.Range("E29").Value = "  +3.34%  "
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$origStyle_D31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.07"
$ws.Range("D31").Style = $origStyle_D31
$ws.Range("E31").Value = "  +26.62%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$origStyle_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.51"
$ws.Range("D32").Style = $origStyle_D32
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$origStyle_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0614"
$ws.Range("D33").Style = $origStyle_D33
$ws.Range("E33").Value = "  +3.92%  "
$ws.Range("B34").Value = "Gas"
$ws.Range("C34").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$origStyle_D34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.34"
$ws.Range("D34").Style = $origStyle_D34
$ws.Range("E34").Value = "  -4.79%  "
$origStyle_D35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0896"
$ws.Range("D35").Style = $origStyle_D35
$ws.Range("E35").Value = "  +10.60%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$origStyle_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.32"
$ws.Range("D37").Style = $origStyle_D37
$ws.Range("E37").Value = "  +9.06%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$origStyle_D38 = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.85"
$ws.Range("D38").Style = $origStyle_D38
$ws.Range("E38").Value = "  -0.54%  "
$origStyle_D39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.12"
$ws.Range("D39").Style = $origStyle_D39
$ws.Range("E39").Value = "  +2.84%  "
$ws.Range("E40").Value = "  +0.54%  "
$origStyle_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.46"
$ws.Range("D41").Style = $origStyle_D41
$ws.Range("E41").Value = "  +15.26%  "
$origStyle_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0226"
$ws.Range("D42").Style = $origStyle_D42
$ws.Range("E42").Value = "  +4.31%  "
$origStyle_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.17"
$ws.Range("D43").Style = $origStyle_D43
$ws.Range("E43").Value = "  +5.25%  "
$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.12"
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = "  +3.20%  "
$origStyle_D45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0910"
$ws.Range("D45").Style = $origStyle_D45
$ws.Range("E45").Value = "  +11.67%  "
$origStyle_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.77"
$ws.Range("D46").Style = $origStyle_D46
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("D47").Value = "1.323.42"
$ws.Range("E47").Value = "  +0.22%  "
$origStyle_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.98"
$ws.Range("D48").Style = $origStyle_D48
$ws.Range("E48").Value = "  +7.50%  "
$origStyle_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.86"
$ws.Range("D49").Style = $origStyle_D49
$ws.Range("E49").Value = "  +85.63%  "
$origStyle_D50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.07"
$ws.Range("D50").Style = $origStyle_D50
$ws.Range("E50").Value = "  +14.48%  "
$ws.Range("D51").Value = "2.294.74"
$ws.Range("E51").Value = "  +4.46%  "
